$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.687213540077209
$ws.Range("B1").Value = 2.06021523475647
$ws.Range("C1").Value = 5.147397518157959
$ws.Range("D1").Value = 1.378853440284729
$ws.Range("E1").Value = 0.6622394323348999
